# Antibody Results test workbook update
# - Header F1: "Data Links" -> "Data Links (but used for Comments here)"
# - Column F (rows 2-10): replace the old "20200617_1" link placeholder with
#   descriptive per-row test comments
# - Add a new row 11: SpecimenAntibodyResults8 / G814450908, testing that a
#   mismatched uploaded Well Position (H10, highlighted like the other
#   "expect error" rows) throws an error
# - Widen column F to fit the new, longer comment text
# - Move the active selection to B12

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 11: mismatched Well Position test case ---
# Copy the formatting from row 9 (same highlight pattern: only the
# "Source Well" cell is flagged) before filling in the new values.
$ws.Range("A9:I9").Copy() | Out-Null
$ws.Range("A11:I11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Specimen ID first (placeholder entered before the comment rewrite pass)
$ws.Range("A11").Value = "SpecimenAntibodyResults8"

# --- Header ---
$ws.Range("F1").Value = "Data Links (but used for Comments here)"

# --- Updated per-row comments (column F) ---
$ws.Range("F2").Value  = "Verify Negative conclusion"
$ws.Range("F3").Value  = "Verify Positive conclusion"
$ws.Range("F4").Value  = "Verify Non-Negative conclusion"
$ws.Range("F6").Value  = "Verify error on missing Specimen ID"
$ws.Range("F7").Value  = "Verify error on missing Biobank Tube ID"
$ws.Range("F8").Value  = "Verify error on missing Conclusion"
$ws.Range("F9").Value  = "Verify error on missing Well Position"
$ws.Range("F10").Value = "Verify error on missing Well Plate Barcode"
$ws.Range("F11").Value = "Verify error when uploaded result Well Position doesn't match"

# --- Remaining new row 11 values ---
$ws.Range("E11").Value = "H10"
$ws.Range("B11").Value = "G814450908"
$ws.Range("C11").Value = "Negative"
$ws.Range("D11").Value = 1
$ws.Range("G11").Value = "AntibodyResults1"
$ws.Range("H11").Value = "Rack 03"

# I11 keeps the quotePrefix-style formatting used throughout column I
# ("RackPos 05", same as rows 6/7/10); pasting value+format from one of
# those rows avoids clobbering the quote-prefix flag the way a plain
# .Value assignment would.
$ws.Range("I7").Copy() | Out-Null
$ws.Range("I11").PasteSpecial(-4104) | Out-Null
$excel.CutCopyMode = 0

# --- Column F width to fit the longer comment text ---
# (target best-fit width is ~50.16; the engine quantizes ColumnWidth to
# 1/6-character steps, so 49.3 is the input that lands closest to it)
$ws.Columns.Item(6).ColumnWidth = 49.3

# --- Move selection ---
$ws.Range("B12").Select() | Out-Null
